# CIERRE 27 MAY 22
# Update the "REMISIONES MAYO 2022" credit-tracking sheet: mark several
# remisiones as paid (fill in payment date + payment amount so the
# outstanding-balance formula in column H nets to 0), and fill in the
# newly-registered remisiones in rows 27-34 (date, client, amount owed,
# and — where paid — the payment date/amount).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("  REMISIONES   MAYO   2022   ")
$ws.Activate()

# --- Mark previously-outstanding remisiones as paid ---------------------
# Row 4
$ws.Range("F4").Value = 44702
$ws.Range("G4").Value = 303

# Row 18
$ws.Range("F18").Value = 44702
$ws.Range("G18").Value = 464

# Row 23
$ws.Range("F23").Value = 44700
$ws.Range("G23").Value = 11923

# Row 24
$ws.Range("F24").Value = 44702
$ws.Range("G24").Value = 410

# Row 26
$ws.Range("F26").Value = 44702
$ws.Range("G26").Value = 4445

# --- Fill in the newly-registered remisiones (rows 27-34) ---------------
# Row 27
$ws.Range("A27").Value = 44700
$ws.Range("D27").Value = "HERRADURA GUSTAVO"
$ws.Range("E27").Value = 12050
$ws.Range("F27").Value = 44701
$ws.Range("G27").Value = 12050

# Row 28
$ws.Range("A28").Value = 44700
$ws.Range("D28").Value = "HERRADURA GUSTAVO"
$ws.Range("E28").Value = 6615
$ws.Range("F28").Value = 44701
$ws.Range("G28").Value = 6615

# Row 29
$ws.Range("A29").Value = 44701
$ws.Range("D29").Value = "HERRADURA GUSTAVO"
$ws.Range("E29").Value = 25704
$ws.Range("F29").Value = 44702
$ws.Range("G29").Value = 25704

# Row 30
$ws.Range("A30").Value = 44702
$ws.Range("D30").Value = "HERRADURA GUSTAVO"
$ws.Range("E30").Value = 14967
$ws.Range("F30").Value = 44703
$ws.Range("G30").Value = 14967

# Row 31
$ws.Range("A31").Value = 44702
$ws.Range("D31").Value = "MAURO"
$ws.Range("E31").Value = 2579
$ws.Range("F31").Value = 44703
$ws.Range("G31").Value = 2579

# Row 32
$ws.Range("A32").Value = 44703
$ws.Range("D32").Value = "GABRIEL"
$ws.Range("E32").Value = 5922
$ws.Range("F32").Value = 44704
$ws.Range("G32").Value = 5922

# Row 33 - still unpaid, only the remision itself is registered
$ws.Range("A33").Value = 44703
$ws.Range("D33").Value = "HERRADURA GUSTAVO"
$ws.Range("E33").Value = 14579

# Row 34
$ws.Range("A34").Value = 44704
$ws.Range("D34").Value = "HERRADURA GUSTAVO"
$ws.Range("E34").Value = 8482
$ws.Range("F34").Value = 44704
$ws.Range("G34").Value = 8482

# --- Restore the view/selection state captured at closing time ----------
$ws.Application.ActiveWindow.ScrollRow = 13
$ws.Range("G35").Select()
